$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: the paragraph holding the 14fig04.jpg figure (right after the
# paragraph ending "...多层级皮肤事实上是以不同的方式吸收和散射光的，就如下图所示：")
# loses the stray paragraph-mark run-properties (<w:rPr><w:rFonts
# w:hint="eastAsia"/></w:rPr>) that used to sit inside its <w:pPr>.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute("进一步复杂化这个过程", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $anchor.Start -and $candidate.Range.End -ge $anchor.End) {
        $anchorIndex = $i
        break
    }
}

$figurePara = $d.Paragraphs.Item($anchorIndex).Next()
$figureRange = $d.Range($figurePara.Range.Start, $figurePara.Range.End)

$figureXml = '<w:p w:rsidR="00D6382D" w:rsidRDefault="00D6382D" w:rsidP="00A818AC"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="840"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="2860675" cy="1752600"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="3" name="图片 3" descr="14fig04.jpg"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1" descr="14fig04.jpg"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="2860675" cy="1752600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$figurePackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:body>' + $figureXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $figureRange.InsertXML($figurePackage)

# ---------------------------------------------------------------------------
# Hunk 2: the closing paragraph ("因为BRDF的定义，" + the _GoBack bookmark) is
# expanded with the rest of the "Rendering with a BRDF" discussion, a new
# "Fresnel Reflectance for Rendering Skin" heading, and its own discussion
# paragraph (which keeps the trailing bookmark).
# ---------------------------------------------------------------------------
$tailPara = $d.Paragraphs.Last
$tailRange = $d.Range($tailPara.Range.Start, $tailPara.Range.End)

$tailXml = '<w:p w:rsidR="00397A71" w:rsidRPr="00397A71" w:rsidRDefault="001440D1" w:rsidP="00397A71"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="840"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>因为BRDF的定义，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>除了计算BRDF本身外，还需要计算点乘（N，L[</w:t></w:r><w:r><w:t>i]</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>）项。同时，还会给每个光源增加一个距离</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>衰减项去根据</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>与光源的距离降低入射光的强度。这个代码同时适用于可以计算L向量和阴影的点光源、平行光源和聚光灯。</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>粗糙表面的镜面反射是由环境贴图光源或区域光反射造成的，这种反射非常复杂，计算起来非常昂贵，我们没有把它们纳入我们的皮肤渲染系统。</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="4"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="180" w:after="180"/><w:ind w:left="420" w:firstLine="420"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="004E49"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="004E49"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Fresnel Reflectance for Rendering Skin</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="840"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>所有基于物理的specular</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>BRDF模型包含一个Fresnel项，通常不会进行详细的解释。这应该是一个非偏振的</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>介电质</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>菲涅尔反射函数，其F0参数为0.028。这来自于比尔定律，且假设皮肤的折射率为1.4。</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$tailPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $tailXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $tailRange.InsertXML($tailPackage)

Write-Host "OK"
